# Update cryptocurrency Price (D) and Volume(1h) (E) columns per the data refresh.
# Values are forced to remain plain text (matching the original inlineStr cells);
# a leading apostrophe is used only where the new text would otherwise be
# auto-parsed by Excel as a number/percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.314.34"
$ws.Range("E2").Value = "  -3.32%  "
$ws.Range("D3").Value = "1.929.72"
$ws.Range("E3").Value = "  -3.87%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'248.67"
$ws.Range("D6").Value = "'0.7202"
$ws.Range("E6").Value = "  -8.40%  "
$ws.Range("D8").Value = "'0.3302"
$ws.Range("D9").Value = "'28.03"
$ws.Range("E9").Value = "  -2.56%  "
$ws.Range("D10").Value = "'0.06923"
$ws.Range("E10").Value = "  -2.47%  "
$ws.Range("D11").Value = "'0.8030"
$ws.Range("E11").Value = "  -6.41%  "
$ws.Range("D12").Value = "'0.08078"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("D13").Value = "1.929.56"
$ws.Range("E13").Value = "  -3.81%  "
$ws.Range("D14").Value = "'5.412"
$ws.Range("E14").Value = "  -3.70%  "
$ws.Range("D15").Value = "'94.73"
$ws.Range("E15").Value = "  -6.29%  "
$ws.Range("D16").Value = "'14.51"
$ws.Range("E16").Value = "  -2.23%  "
$ws.Range("D17").Value = "30.298.95"
$ws.Range("D18").Value = "'0.000008292"
$ws.Range("E18").Value = "  +4.27%  "
$ws.Range("D19").Value = "'252.73"
$ws.Range("E19").Value = "  -8.54%  "
$ws.Range("E20").Value = "  -2.42%  "
$ws.Range("D21").Value = "2.183.63"
$ws.Range("E21").Value = "  -3.93%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'6.868"
$ws.Range("D25").Value = "'9.722"
$ws.Range("E25").Value = "  -3.67%  "
$ws.Range("D26").Value = "'159.56"
$ws.Range("E26").Value = "  -3.01%  "
$ws.Range("D27").Value = "'2.392"
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("D28").Value = "'19.15"
$ws.Range("E28").Value = "  -4.39%  "
$ws.Range("D29").Value = "'0.1330"
$ws.Range("E29").Value = "  -12.25%  "
$ws.Range("D30").Value = "'1.553"
$ws.Range("E30").Value = "  -4.67%  "
$ws.Range("D31").Value = "'1.338"
$ws.Range("E31").Value = "  -1.60%  "
$ws.Range("D32").Value = "'4.402"
$ws.Range("E32").Value = "  -5.20%  "
$ws.Range("D33").Value = "'4.195"
$ws.Range("E33").Value = "  -5.21%  "
$ws.Range("D34").Value = "'0.05111"
$ws.Range("E34").Value = "  -2.83%  "
$ws.Range("D35").Value = "'1.219"
$ws.Range("E35").Value = "  -0.43%  "
$ws.Range("E36").Value = "  -4.38%  "
$ws.Range("D37").Value = "'2.742"
$ws.Range("E37").Value = "  -2.34%  "
$ws.Range("D38").Value = "'0.01977"
$ws.Range("E38").Value = "  -2.17%  "
$ws.Range("D39").Value = "'2.830"
$ws.Range("E39").Value = "  -3.88%  "
$ws.Range("D40").Value = "'6.590"
$ws.Range("E40").Value = "  -2.42%  "
$ws.Range("D41").Value = "'78.85"
$ws.Range("E41").Value = "  -2.34%  "
$ws.Range("D42").Value = "'0.4465"
$ws.Range("E42").Value = "  -6.33%  "
$ws.Range("D43").Value = "'1.992"
$ws.Range("E43").Value = "  -8.44%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "'0.8363"
$ws.Range("E45").Value = "  -2.64%  "
$ws.Range("D46").Value = "'102.06"
$ws.Range("E46").Value = "  -3.85%  "
$ws.Range("D47").Value = "'9.755"
$ws.Range("E47").Value = "  -1.93%  "
$ws.Range("D48").Value = "'7.296"
$ws.Range("E48").Value = "  -6.30%  "
$ws.Range("D49").Value = "'36.61"
$ws.Range("E49").Value = "  -1.30%  "
$ws.Range("D50").Value = "'0.05955"
$ws.Range("D51").Value = "'0.4085"
$ws.Range("E51").Value = "  -6.85%  "
